$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B138").Value = "s137_e105_24019-32_2_7.jpeg"
$ws.Range("D138").Value = "'2532"
$ws.Range("E138").Value = "'1403"
$ws.Range("H138").Value = "'139"
$ws.Range("B139").Value = "s138_e111_24019-32_2_1.jpeg"
$ws.Range("D139").Value = "'2346"
$ws.Range("E139").Value = "'132"
$ws.Range("H139").Value = "'156"
$ws.Range("B140").Value = "s139_e104_24019-32_1_8.jpeg"
$ws.Range("D140").Value = "'701"
$ws.Range("E140").Value = "'1350"
$ws.Range("H140").Value = "'138"
$ws.Range("B141").Value = "s140_e109_24019-32_2_3.jpeg"
$ws.Range("D141").Value = "'2525"
$ws.Range("E141").Value = "'1061"
$ws.Range("H141").Value = "'61"
$ws.Range("B142").Value = "s141_e119_24019-32_3_7.jpeg"
$ws.Range("D142").Value = "'1927"
$ws.Range("E142").Value = "'1511"
$ws.Range("H142").Value = "'78"
$ws.Range("B143").Value = "s142_e118_24019-32_3_6.jpeg"
$ws.Range("D143").Value = "'994"
$ws.Range("E143").Value = "'826"
$ws.Range("H143").Value = "'105"
$ws.Range("B144").Value = "s143_e124_24019-32_4_4.jpeg"
$ws.Range("D144").Value = "'481"
$ws.Range("E144").Value = "'1494"
$ws.Range("H144").Value = "'63"
$ws.Range("B145").Value = "s144_e122_24019-32_4_6.jpeg"
$ws.Range("D145").Value = "'119"
$ws.Range("E145").Value = "'1109"
$ws.Range("H145").Value = "'57"
$ws.Range("B146").Value = "s145_e113_24019-32_3_1.jpeg"
$ws.Range("D146").Value = "'513"
$ws.Range("E146").Value = "'197"
$ws.Range("H146").Value = "'143"
$ws.Range("B147").Value = "s146_e98_24019-32_1_2.jpeg"
$ws.Range("D147").Value = "'536"
$ws.Range("E147").Value = "'1753"
$ws.Range("H147").Value = "'46"
$ws.Range("D148").Value = "'114"
$ws.Range("E148").Value = "'241"
$ws.Range("H148").Value = "'154"
$ws.Range("B149").Value = "s148_e101_24019-32_1_5.jpeg"
$ws.Range("D149").Value = "'2207"
$ws.Range("E149").Value = "'1159"
$ws.Range("H149").Value = "'111"
$ws.Range("B150").Value = "s149_e112_24019-32_2_0.jpeg"
$ws.Range("D150").Value = "'1354"
$ws.Range("E150").Value = "'1185"
$ws.Range("H150").Value = "'138"
$ws.Range("B151").Value = "s150_e117_24019-32_3_5.jpeg"
$ws.Range("D151").Value = "'2358"
$ws.Range("E151").Value = "'1436"
$ws.Range("H151").Value = "'53"
$ws.Range("B152").Value = "s151_e116_24019-32_3_4.jpeg"
$ws.Range("D152").Value = "'1148"
$ws.Range("E152").Value = "'1498"
$ws.Range("H152").Value = "'146"
$ws.Range("B153").Value = "s152_e107_24019-32_2_5.jpeg"
$ws.Range("D153").Value = "'2506"
$ws.Range("E153").Value = "'225"
$ws.Range("H153").Value = "'74"
$ws.Range("B154").Value = "s153_e142_24019-32_3_2.jpeg"
$ws.Range("D154").Value = "'464"
$ws.Range("E154").Value = "'661"
$ws.Range("H154").Value = "'78"
$ws.Range("B155").Value = "s154_e146_24019-32_3_6.jpeg"
$ws.Range("D155").Value = "'312"
$ws.Range("E155").Value = "'246"
$ws.Range("H155").Value = "'169"
$ws.Range("B156").Value = "s155_e150_24019-32_4_6.jpeg"
$ws.Range("D156").Value = "'1022"
$ws.Range("E156").Value = "'413"
$ws.Range("H156").Value = "'114"
$ws.Range("B157").Value = "s156_e144_24019-32_3_4.jpeg"
$ws.Range("D157").Value = "'1839"
$ws.Range("E157").Value = "'307"
$ws.Range("H157").Value = "'152"
$ws.Range("B158").Value = "s157_e125_24019-32_1_1.jpeg"
$ws.Range("D158").Value = "'590"
$ws.Range("E158").Value = "'727"
$ws.Range("H158").Value = "'79"
$ws.Range("B159").Value = "s158_e130_24019-32_1_6.jpeg"
$ws.Range("D159").Value = "'302"
$ws.Range("E159").Value = "'1355"
$ws.Range("H159").Value = "'19"
$ws.Range("B160").Value = "s159_e127_24019-32_1_3.jpeg"
$ws.Range("D160").Value = "'1954"
$ws.Range("E160").Value = "'1233"
$ws.Range("H160").Value = "'143"
$ws.Range("B161").Value = "s160_e132_24019-32_1_8.jpeg"
$ws.Range("D161").Value = "'1363"
$ws.Range("E161").Value = "'913"
$ws.Range("H161").Value = "'26"
$ws.Range("B162").Value = "s161_e135_24019-32_2_5.jpeg"
$ws.Range("D162").Value = "'1045"
$ws.Range("E162").Value = "'744"
$ws.Range("H162").Value = "'53"
$ws.Range("B163").Value = "s162_e145_24019-32_3_5.jpeg"
$ws.Range("D163").Value = "'402"
$ws.Range("E163").Value = "'1109"
$ws.Range("H163").Value = "'164"
$ws.Range("B164").Value = "s163_e133_24019-32_2_7.jpeg"
$ws.Range("D164").Value = "'445"
$ws.Range("E164").Value = "'346"
$ws.Range("H164").Value = "'49"
$ws.Range("B165").Value = "s164_e137_24019-32_2_3.jpeg"
$ws.Range("D165").Value = "'1098"
$ws.Range("E165").Value = "'95"
$ws.Range("H165").Value = "'101"
$ws.Range("B166").Value = "s165_e129_24019-32_1_5.jpeg"
$ws.Range("D166").Value = "'378"
$ws.Range("E166").Value = "'988"
$ws.Range("H166").Value = "'172"
$ws.Range("B167").Value = "s166_e134_24019-32_2_6.jpeg"
$ws.Range("D167").Value = "'1044"
$ws.Range("E167").Value = "'484"
$ws.Range("H167").Value = "'84"
$ws.Range("B168").Value = "s167_e128_24019-32_1_4.jpeg"
$ws.Range("D168").Value = "'1365"
$ws.Range("E168").Value = "'993"
$ws.Range("H168").Value = "'78"
$ws.Range("B169").Value = "s168_e131_24019-32_1_7.jpeg"
$ws.Range("D169").Value = "'386"
$ws.Range("E169").Value = "'628"
$ws.Range("H169").Value = "'92"

"done"